$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the whole "Caption" bullet paragraph (its content is
#    dropped entirely, leaving the preceding "Geolocation" bullet as
#    the last one before the "Insights from the joined dataset"
#    heading).
# ------------------------------------------------------------------
$capRng = $d.Content
$foundCap = $capRng.Find.Execute("Caption: - Some images also had metadata fields like", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundCap) {
    $capRng.Expand(4) | Out-Null   # wdParagraph - grab the whole paragraph incl. mark
    $capRng.Delete()
}

# ------------------------------------------------------------------
# 2. Replace the placeholder line "Through Image Captioning ----"
#    with the real answer text.
# ------------------------------------------------------------------
$newText1 = "For UFO stalker dataset , there was no straightforward question which we could think of that this dataset helped to answer."
$d.Content.Find.Execute("Through Image Captioning ----", $true, $false, $false, $false, $false, $true, 1, $false, $newText1, 2) | Out-Null

# Relocate the "_GoBack" bookmark so it sits right after the text we
# just typed (mirrors where Word drops it after the last edit made).
# NOTE: adding a bookmark to a zero-length range that sits exactly at
# a paragraph's end position mis-places it in this host, so we nudge
# past the boundary with a throwaway character, bookmark there, then
# remove the throwaway character.
try {
    $oldBm = $d.Bookmarks("_GoBack")
    $oldBm.Delete()
} catch {
}

$bmRng = $d.Content
$bmRng.Find.Execute($newText1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRng.Collapse(0) | Out-Null               # wdCollapseEnd
$bmRng.InsertAfter("X") | Out-Null
$bmPos = $d.Range($bmRng.Start, $bmRng.Start)
$d.Bookmarks.Add("_GoBack", $bmPos) | Out-Null
$placeholder = $d.Range($bmRng.Start, $bmRng.Start + 1)
$placeholder.Delete()

# ------------------------------------------------------------------
# 3. Replace the other placeholder line "Image Captioning ----"
#    (the OCR/Image-captioning Q&A bullet) with the real answer text.
# ------------------------------------------------------------------
$newText2 = "For Image captioning task, use of docker made it quite easy, but as discussed above, a lot of captions generated were not accurate as they were not able to detect/recognized the sighted UFO. "
$d.Content.Find.Execute("Image Captioning ----", $true, $false, $false, $false, $false, $true, 1, $false, $newText2, 2) | Out-Null

Write-Output "done"
